$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed naive component forecaster bug - Presentation state 11.02.
# Update the per-quarter naive-average QoQ error values (rows 24-52, columns B:K).

# Row 24
$ws.Range("K24").Value = -0.35579190771253

# Row 25
$ws.Range("J25").Value = -1.554235363265292
$ws.Range("K25").Value = -0.2643100270664007

# Row 26
$ws.Range("I26").Value = -1.679044851730669
$ws.Range("J26").Value = -0.3891195155317774
$ws.Range("K26").Value = 1.820076036519961

# Row 27
$ws.Range("H27").Value = -1.398238077646767
$ws.Range("I27").Value = -0.1083127414478752
$ws.Range("J27").Value = 2.100882810603863
$ws.Range("K27").Value = -0.8049070241509069

# Row 28
$ws.Range("G28").Value = -2.979044851730669
$ws.Range("H28").Value = -1.689119515531778
$ws.Range("I28").Value = 0.5200760365199608
$ws.Range("J28").Value = -2.385713798234809
$ws.Range("K28").Value = -2.026036591042376

# Row 29
$ws.Range("F29").Value = -0.9284238248828705
$ws.Range("G29").Value = 0.3615015113160212
$ws.Range("H29").Value = 2.570697063367759
$ws.Range("I29").Value = -0.3350927713870104
$ws.Range("J29").Value = 0.0245844358054228
$ws.Range("K29").Value = -0.3732558837842244

# Row 30
$ws.Range("E30").Value = -0.8773751373745151
$ws.Range("F30").Value = 0.4125501988243765
$ws.Range("G30").Value = 2.621745750876115
$ws.Range("H30").Value = -0.2840440838786551
$ws.Range("I30").Value = 0.07563312331377814
$ws.Range("J30").Value = -0.322207196275869
$ws.Range("K30").Value = 3.13844758671037

# Row 31
$ws.Range("D31").Value = -0.828109867557302
$ws.Range("E31").Value = 0.4618154686415896
$ws.Range("F31").Value = 2.671011020693328
$ws.Range("G31").Value = -0.234778814061442
$ws.Range("H31").Value = 0.1248983931309913
$ws.Range("I31").Value = -0.2729419264586559
$ws.Range("J31").Value = 3.187712856527583
$ws.Range("K31").Value = 0.257529852677735

# Row 32
$ws.Range("C32").Value = -2.179044851730669
$ws.Range("D32").Value = -0.8891195155317775
$ws.Range("E32").Value = 1.320076036519961
$ws.Range("F32").Value = -1.585713798234809
$ws.Range("G32").Value = -1.226036591042376
$ws.Range("H32").Value = -1.623876910632023
$ws.Range("I32").Value = 1.836777872354216
$ws.Range("J32").Value = -1.093405131495632
$ws.Range("K32").Value = -1.329975481959951

# Row 33
$ws.Range("B33").Value = -1.292712606546062
$ws.Range("C33").Value = -0.002787270347170079
$ws.Range("D33").Value = 2.206408281704568
$ws.Range("E33").Value = -0.6993815530502017
$ws.Range("F33").Value = -0.3397043458577684
$ws.Range("G33").Value = -0.7375446654474156
$ws.Range("H33").Value = 2.723110117538823
$ws.Range("I33").Value = -0.2070728863110247
$ws.Range("J33").Value = -0.4436432367753436
$ws.Range("K33").Value = -0.1180488047942703

# Row 34
$ws.Range("B34").Value = 1.289925336198892
$ws.Range("C34").Value = 3.49912088825063
$ws.Range("D34").Value = 0.59333105349586
$ws.Range("E34").Value = 0.9530082606882933
$ws.Range("F34").Value = 0.5551679410986461
$ws.Range("G34").Value = 4.015822724084884
$ws.Range("H34").Value = 1.085639720235037
$ws.Range("I34").Value = 0.8490693697707181
$ws.Range("J34").Value = 1.174663801751791
$ws.Range("K34").Value = 1.738389274431257

# Row 35
$ws.Range("B35").Value = 2.209195552051738
$ws.Range("C35").Value = -0.6965942827030316
$ws.Range("D35").Value = -0.3369170755105984
$ws.Range("E35").Value = -0.7347573951002455
$ws.Range("F35").Value = 2.725897387885993
$ws.Range("G35").Value = -0.2042856159638546
$ws.Range("H35").Value = -0.4408559664281735
$ws.Range("I35").Value = -0.1152615344471002
$ws.Range("J35").Value = 0.448463938232365
$ws.Range("K35").Value = 0.162308343228304

# Row 36
$ws.Range("B36").Value = -2.90578983475477
$ws.Range("C36").Value = -2.546112627562337
$ws.Range("D36").Value = -2.943952947151984
$ws.Range("E36").Value = 0.516701835834255
$ws.Range("F36").Value = -2.413481168015593
$ws.Range("G36").Value = -2.650051518479912
$ws.Range("H36").Value = -2.324457086498839
$ws.Range("I36").Value = -1.760731613819373
$ws.Range("J36").Value = -2.046887208823434
$ws.Range("K36").Value = -1.925041643302435

# Row 37
$ws.Range("B37").Value = 0.3596772071924332
$ws.Range("C37").Value = -0.03816311239721393
$ws.Range("D37").Value = 3.422491670589025
$ws.Range("E37").Value = 0.4923086667391769
$ws.Range("F37").Value = 0.2557383162748581
$ws.Range("G37").Value = 0.5813327482559314
$ws.Range("H37").Value = 1.145058220935397
$ws.Range("I37").Value = 0.8589026259313355
$ws.Range("J37").Value = 0.9807481914523348
$ws.Range("K37").Value = 0.5470686472140469

# Row 38
$ws.Range("B38").Value = -0.3978403195896472
$ws.Range("C38").Value = 3.062814463396592
$ws.Range("D38").Value = 0.1326314595467437
$ws.Range("E38").Value = -0.1039388909175751
$ws.Range("F38").Value = 0.2216555410634982
$ws.Range("G38").Value = 0.7853810137429633
$ws.Range("H38").Value = 0.4992254187389023
$ws.Range("I38").Value = 0.6210709842599016
$ws.Range("J38").Value = 0.1873914400216137
$ws.Range("K38").Value = 2.073267678908721

# Row 39
$ws.Range("B39").Value = 3.460654782986239
$ws.Range("C39").Value = 0.5304717791363909
$ws.Range("D39").Value = 0.293901428672072
$ws.Range("E39").Value = 0.6194958606531453
$ws.Range("F39").Value = 1.18322133333261
$ws.Range("G39").Value = 0.8970657383285494
$ws.Range("H39").Value = 1.018911303849549
$ws.Range("I39").Value = 0.5852317596112608
$ws.Range("J39").Value = 2.471107998498368
$ws.Range("K39").Value = 0.7675489330019185

# Row 40
$ws.Range("B40").Value = -2.930183003849848
$ws.Range("C40").Value = -3.166753354314167
$ws.Range("D40").Value = -2.841158922333094
$ws.Range("E40").Value = -2.277433449653628
$ws.Range("F40").Value = -2.563589044657689
$ws.Range("G40").Value = -2.44174347913669
$ws.Range("H40").Value = -2.875423023374978
$ws.Range("I40").Value = -0.9895467844878709
$ws.Range("J40").Value = -2.69310584998432
$ws.Range("K40").Value = -2.748687381508546

# Row 41
$ws.Range("B41").Value = -0.2365703504643188
$ws.Range("C41").Value = 0.08902408151675445
$ws.Range("D41").Value = 0.6527495541962196
$ws.Range("E41").Value = 0.3665939591921586
$ws.Range("F41").Value = 0.4884395247131579
$ws.Range("G41").Value = 0.05475998047486996
$ws.Range("H41").Value = 1.940636219361977
$ws.Range("I41").Value = 0.2370771538655276
$ws.Range("J41").Value = 0.181495622341302
$ws.Range("K41").Value = 0.9684454849153923

# Row 42
$ws.Range("B42").Value = 0.3255944319810733
$ws.Range("C42").Value = 0.8893199046605385
$ws.Range("D42").Value = 0.6031643096564774
$ws.Range("E42").Value = 0.7250098751774767
$ws.Range("F42").Value = 0.2913303309391888
$ws.Range("G42").Value = 2.177206569826296
$ws.Range("H42").Value = 0.4736475043298465
$ws.Range("I42").Value = 0.4180659728056209
$ws.Range("J42").Value = 1.205015835379711
$ws.Range("K42").Value = -0.9846815754178531

# Row 43
$ws.Range("B43").Value = 0.5637254726794652
$ws.Range("C43").Value = 0.2775698776754041
$ws.Range("D43").Value = 0.3994154431964034
$ws.Range("E43").Value = -0.03426410104188449
$ws.Range("F43").Value = 1.851612137845223
$ws.Range("G43").Value = 0.1480530723487732
$ws.Range("H43").Value = 0.09247154082454756
$ws.Range("I43").Value = 0.8794214033986378
$ws.Range("J43").Value = -1.310276007398926
$ws.Range("K43").Value = -0.4722849543269805

# Row 44
$ws.Range("B44").Value = -0.286155595004061
$ws.Range("C44").Value = -0.1643100294830617
$ws.Range("D44").Value = -0.5979895737213496
$ws.Range("E44").Value = 1.287886665165757
$ws.Range("F44").Value = -0.415672400330692
$ws.Range("G44").Value = -0.4712539318549176
$ws.Range("H44").Value = 0.3156959307191727
$ws.Range("I44").Value = -1.874001480078392
$ws.Range("J44").Value = -1.036010427006446

# Row 45
$ws.Range("B45").Value = 0.1218455655209993
$ws.Range("C45").Value = -0.3118339787172886
$ws.Range("D45").Value = 1.574042260169819
$ws.Range("E45").Value = -0.129516805326631
$ws.Range("F45").Value = -0.1850983368508566
$ws.Range("G45").Value = 0.6018515257232337
$ws.Range("H45").Value = -1.58784588507433
$ws.Range("I45").Value = -0.7498548320023846

# Row 46
$ws.Range("B46").Value = -0.4336795442382879
$ws.Range("C46").Value = 1.452196694648819
$ws.Range("D46").Value = -0.2513623708476302
$ws.Range("E46").Value = -0.3069439023718559
$ws.Range("F46").Value = 0.4800059602022344
$ws.Range("G46").Value = -1.70969145059533
$ws.Range("H46").Value = -0.8717003975233839

# Row 47
$ws.Range("B47").Value = 1.885876238887107
$ws.Range("C47").Value = 0.1823171733906577
$ws.Range("D47").Value = 0.126735641866432
$ws.Range("E47").Value = 0.9136855044405223
$ws.Range("F47").Value = -1.276011906357042
$ws.Range("G47").Value = -0.438020853285096

# Row 48
$ws.Range("B48").Value = -1.703559065496449
$ws.Range("C48").Value = -1.759140597020675
$ws.Range("D48").Value = -0.9721907344465848
$ws.Range("E48").Value = -3.161888145244149
$ws.Range("F48").Value = -2.323897092172203

# Row 49
$ws.Range("B49").Value = -0.05558153152422562
$ws.Range("C49").Value = 0.7313683310498647
$ws.Range("D49").Value = -1.4583290797477
$ws.Range("E49").Value = -0.6203380266757537

# Row 50
$ws.Range("B50").Value = 0.7869498625740903
$ws.Range("C50").Value = -1.402747548223474
$ws.Range("D50").Value = -0.564756495151528

# Row 51
$ws.Range("B51").Value = -2.189697410797564
$ws.Range("C51").Value = -1.351706357725618

# Row 52
$ws.Range("B52").Value = 0.8379910530719459

